# Update NATMI LR-pair stats (Ptprz1-L1cam) with recalculated TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.01848533333333334
$ws.Range("H2").Value = 0.05545600000000001
$ws.Range("I2").Value = 0.001625201930372746
$ws.Range("J2").Value = 0.001625201930372746
$ws.Range("M2").Value = 3.685507
$ws.Range("N2").Value = 11.056521
$ws.Range("O2").Value = 0.3585631737883472
$ws.Range("P2").Value = 0.3585631737883472
$ws.Range("Q2").Value = 0.06812782539733334
$ws.Range("R2").Value = 0.6131504285760001
$ws.Range("S2").Value = 0.0005827375622014003
$ws.Range("T2").Value = 0.0005827375622014002
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.01848533333333334
$ws.Range("H3").Value = 0.05545600000000001
$ws.Range("I3").Value = 0.001625201930372746
$ws.Range("J3").Value = 0.001625201930372746
$ws.Range("O3").Value = 0.00964718443071163
$ws.Range("P3").Value = 0.00964718443071163
$ws.Range("Q3").Value = 0.001832987168
$ws.Range("R3").Value = 0.016496884512
$ws.Range("S3").Value = 0.00001567862275945444
$ws.Range("T3").Value = 0.00001567862275945444
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.01848533333333334
$ws.Range("H4").Value = 0.05545600000000001
$ws.Range("I4").Value = 0.001625201930372746
$ws.Range("J4").Value = 0.001625201930372746
$ws.Range("M4").Value = 6.493877
$ws.Range("N4").Value = 19.481631
$ws.Range("O4").Value = 0.6317896417809412
$ws.Range("P4").Value = 0.6317896417809411
$ws.Range("Q4").Value = 0.1200414809706667
$ws.Range("R4").Value = 1.080373328736
$ws.Range("S4").Value = 0.001026785745411891
$ws.Range("T4").Value = 0.001026785745411891
$ws.Range("I5").Value = 0.002698334581238102
$ws.Range("J5").Value = 0.002698334581238102
$ws.Range("M5").Value = 3.685507
$ws.Range("N5").Value = 11.056521
$ws.Range("O5").Value = 0.3585631737883472
$ws.Range("P5").Value = 0.3585631737883472
$ws.Range("Q5").Value = 0.1131131238393333
$ws.Range("R5").Value = 1.018018114554
$ws.Range("S5").Value = 0.0009675234113915847
$ws.Range("T5").Value = 0.0009675234113915846
$ws.Range("I6").Value = 0.002698334581238102
$ws.Range("J6").Value = 0.002698334581238102
$ws.Range("O6").Value = 0.00964718443071163
$ws.Range("P6").Value = 0.00964718443071163
$ws.Range("S6").Value = 0.000026031331360971
$ws.Range("T6").Value = 0.000026031331360971
$ws.Range("I7").Value = 0.002698334581238102
$ws.Range("J7").Value = 0.002698334581238102
$ws.Range("M7").Value = 6.493877
$ws.Range("N7").Value = 19.481631
$ws.Range("O7").Value = 0.6317896417809412
$ws.Range("P7").Value = 0.6317896417809411
$ws.Range("Q7").Value = 0.1993057436326667
$ws.Range("R7").Value = 1.793751692694
$ws.Range("S7").Value = 0.001704779838485546
$ws.Range("T7").Value = 0.001704779838485546
$ws.Range("G8").Value = 11.32499966666667
$ws.Range("H8").Value = 33.974999
$ws.Range("I8").Value = 0.9956764634883892
$ws.Range("J8").Value = 0.995676463488389
$ws.Range("M8").Value = 3.685507
$ws.Range("N8").Value = 11.056521
$ws.Range("O8").Value = 0.3585631737883472
$ws.Range("P8").Value = 0.3585631737883472
$ws.Range("Q8").Value = 41.73836554649768
$ws.Range("R8").Value = 375.645289918479
$ws.Range("S8").Value = 0.3570129128147543
$ws.Range("T8").Value = 0.3570129128147542
$ws.Range("G9").Value = 11.32499966666667
$ws.Range("H9").Value = 33.974999
$ws.Range("I9").Value = 0.9956764634883892
$ws.Range("J9").Value = 0.995676463488389
$ws.Range("O9").Value = 0.00964718443071163
$ws.Range("P9").Value = 0.00964718443071163
$ws.Range("Q9").Value = 1.122975641947
$ws.Range("R9").Value = 10.106780777523
$ws.Range("S9").Value = 0.009605474476591205
$ws.Range("T9").Value = 0.009605474476591202
$ws.Range("G10").Value = 11.32499966666667
$ws.Range("H10").Value = 33.974999
$ws.Range("I10").Value = 0.9956764634883892
$ws.Range("J10").Value = 0.995676463488389
$ws.Range("M10").Value = 6.493877
$ws.Range("N10").Value = 19.481631
$ws.Range("O10").Value = 0.6317896417809412
$ws.Range("P10").Value = 0.6317896417809411
$ws.Range("Q10").Value = 73.54315486037434
$ws.Range("R10").Value = 661.8883937433691
$ws.Range("S10").Value = 0.6290580761970438
$ws.Range("T10").Value = 0.6290580761970436
